$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.168739795684814
$ws.Range("B1").Value = 1.083346486091614
$ws.Range("C1").Value = 15
$ws.Range("D1").Value = 1.997352123260498
$ws.Range("E1").Value = 0.9817475080490112
